$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-07-06 Sunday" "2025-07-07 Monday"

Replace-Text "544×6=" "295×5="
Replace-Text "505×6=" "440×8="
Replace-Text "320×3=" "819×6="
Replace-Text "486×4=" "300×6="
Replace-Text "774×4=" "998×3="
Replace-Text "310×6=" "941×4="
Replace-Text "291×5=" "673×2="
Replace-Text "707×6=" "105×3="
Replace-Text "617×3=" "875×4="
Replace-Text "942×2=" "777×6="
Replace-Text "316×6=" "582×2="
Replace-Text "107×5=" "213×2="
Replace-Text "829×8=" "966×9="
Replace-Text "254×6=" "378×9="
Replace-Text "639×5=" "756×3="
Replace-Text "176×2=" "380×9="
Replace-Text "710×2=" "925×4="
Replace-Text "375×5=" "468×2="
Replace-Text "615×3=" "474×8="
Replace-Text "861×9=" "703×7="
Replace-Text "850×7=" "449×9="
Replace-Text "488×7=" "354×2="
Replace-Text "350×4=" "821×3="
Replace-Text "513×9=" "804×8="
Replace-Text "494×2=" "208×5="
